$d = $word.ActiveDocument

$d.Content.Find.Execute("2022-12-11 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2022-12-12 Monday", 2) | Out-Null
$d.Content.Find.Execute("68-41=", $true, $false, $false, $false, $false, $true, 1, $false, "59-54=", 2) | Out-Null
$d.Content.Find.Execute("97-22=", $true, $false, $false, $false, $false, $true, 1, $false, "42-26=", 2) | Out-Null
$d.Content.Find.Execute("67-38=", $true, $false, $false, $false, $false, $true, 1, $false, "41+0=", 2) | Out-Null
$d.Content.Find.Execute("53+1=", $true, $false, $false, $false, $false, $true, 1, $false, "32+38=", 2) | Out-Null
$d.Content.Find.Execute("69-28=", $true, $false, $false, $false, $false, $true, 1, $false, "65-46=", 2) | Out-Null
$d.Content.Find.Execute("28+70=", $true, $false, $false, $false, $false, $true, 1, $false, "63-27=", 2) | Out-Null
$d.Content.Find.Execute("91-82=", $true, $false, $false, $false, $false, $true, 1, $false, "15+30=", 2) | Out-Null
$d.Content.Find.Execute("12+74=", $true, $false, $false, $false, $false, $true, 1, $false, "48-8=", 2) | Out-Null
$d.Content.Find.Execute("39+3=", $true, $false, $false, $false, $false, $true, 1, $false, "88-10=", 2) | Out-Null
$d.Content.Find.Execute("73-3=", $true, $false, $false, $false, $false, $true, 1, $false, "84-9=", 2) | Out-Null
$d.Content.Find.Execute("86-70=", $true, $false, $false, $false, $false, $true, 1, $false, "5-0=", 2) | Out-Null
$d.Content.Find.Execute("19+37=", $true, $false, $false, $false, $false, $true, 1, $false, "70-39=", 2) | Out-Null
$d.Content.Find.Execute("29+29=", $true, $false, $false, $false, $false, $true, 1, $false, "60+9=", 2) | Out-Null
$d.Content.Find.Execute("75+7=", $true, $false, $false, $false, $false, $true, 1, $false, "22+36=", 2) | Out-Null
$d.Content.Find.Execute("10+17=", $true, $false, $false, $false, $false, $true, 1, $false, "77-11=", 2) | Out-Null
$d.Content.Find.Execute("47+43=", $true, $false, $false, $false, $false, $true, 1, $false, "40+28=", 2) | Out-Null
$d.Content.Find.Execute("87-32=", $true, $false, $false, $false, $false, $true, 1, $false, "50+42=", 2) | Out-Null
$d.Content.Find.Execute("48+51=", $true, $false, $false, $false, $false, $true, 1, $false, "86-82=", 2) | Out-Null
$d.Content.Find.Execute("94-2=", $true, $false, $false, $false, $false, $true, 1, $false, "80-14=", 2) | Out-Null
$d.Content.Find.Execute("16+61=", $true, $false, $false, $false, $false, $true, 1, $false, "50-46=", 2) | Out-Null
$d.Content.Find.Execute("22+12=", $true, $false, $false, $false, $false, $true, 1, $false, "32-28=", 2) | Out-Null
$d.Content.Find.Execute("97-57=", $true, $false, $false, $false, $false, $true, 1, $false, "84-67=", 2) | Out-Null
$d.Content.Find.Execute("63-6=", $true, $false, $false, $false, $false, $true, 1, $false, "16+46=", 2) | Out-Null
$d.Content.Find.Execute("79-53=", $true, $false, $false, $false, $false, $true, 1, $false, "11+84=", 2) | Out-Null
$d.Content.Find.Execute("60+18=", $true, $false, $false, $false, $false, $true, 1, $false, "87-65=", 2) | Out-Null
$d.Content.Find.Execute("51-17=", $true, $false, $false, $false, $false, $true, 1, $false, "44+33=", 2) | Out-Null
$d.Content.Find.Execute("44-15=", $true, $false, $false, $false, $false, $true, 1, $false, "67-10=", 2) | Out-Null
$d.Content.Find.Execute("36-26=", $true, $false, $false, $false, $false, $true, 1, $false, "16+14=", 2) | Out-Null
$d.Content.Find.Execute("58-38=", $true, $false, $false, $false, $false, $true, 1, $false, "6+35=", 2) | Out-Null
$d.Content.Find.Execute("30+30=", $true, $false, $false, $false, $false, $true, 1, $false, "24+70=", 2) | Out-Null
$d.Content.Find.Execute("82-36=", $true, $false, $false, $false, $false, $true, 1, $false, "56-10=", 2) | Out-Null
$d.Content.Find.Execute("61+10=", $true, $false, $false, $false, $false, $true, 1, $false, "52-23=", 2) | Out-Null
$d.Content.Find.Execute("28+25=", $true, $false, $false, $false, $false, $true, 1, $false, "68+30=", 2) | Out-Null
$d.Content.Find.Execute("7+78=", $true, $false, $false, $false, $false, $true, 1, $false, "20+6=", 2) | Out-Null
$d.Content.Find.Execute("62-39=", $true, $false, $false, $false, $false, $true, 1, $false, "35-0=", 2) | Out-Null
$d.Content.Find.Execute("95-87=", $true, $false, $false, $false, $false, $true, 1, $false, "26-12=", 2) | Out-Null
$d.Content.Find.Execute("26+26=", $true, $false, $false, $false, $false, $true, 1, $false, "10+69=", 2) | Out-Null
$d.Content.Find.Execute("56+16=", $true, $false, $false, $false, $false, $true, 1, $false, "26+42=", 2) | Out-Null
$d.Content.Find.Execute("31+68=", $true, $false, $false, $false, $false, $true, 1, $false, "34+36=", 2) | Out-Null
$d.Content.Find.Execute("82-33=", $true, $false, $false, $false, $false, $true, 1, $false, "1+88=", 2) | Out-Null
$d.Content.Find.Execute("50+24=", $true, $false, $false, $false, $false, $true, 1, $false, "25+13=", 2) | Out-Null
$d.Content.Find.Execute("54+39=", $true, $false, $false, $false, $false, $true, 1, $false, "15+13=", 2) | Out-Null
$d.Content.Find.Execute("48+31=", $true, $false, $false, $false, $false, $true, 1, $false, "87-1=", 2) | Out-Null
$d.Content.Find.Execute("39+55=", $true, $false, $false, $false, $false, $true, 1, $false, "39+40=", 2) | Out-Null
$d.Content.Find.Execute("14+44=", $true, $false, $false, $false, $false, $true, 1, $false, "18+41=", 2) | Out-Null
$d.Content.Find.Execute("48-43=", $true, $false, $false, $false, $false, $true, 1, $false, "98-36=", 2) | Out-Null
$d.Content.Find.Execute("98-33=", $true, $false, $false, $false, $false, $true, 1, $false, "23-1=", 2) | Out-Null
$d.Content.Find.Execute("49+23=", $true, $false, $false, $false, $false, $true, 1, $false, "65-20=", 2) | Out-Null
$d.Content.Find.Execute("32+19=", $true, $false, $false, $false, $false, $true, 1, $false, "17+44=", 2) | Out-Null
$d.Content.Find.Execute("18-5=", $true, $false, $false, $false, $false, $true, 1, $false, "37-30=", 2) | Out-Null
$d.Content.Find.Execute("98-13=", $true, $false, $false, $false, $false, $true, 1, $false, "64-24=", 2) | Out-Null
$d.Content.Find.Execute("8+55=", $true, $false, $false, $false, $false, $true, 1, $false, "54+13=", 2) | Out-Null
$d.Content.Find.Execute("34+61=", $true, $false, $false, $false, $false, $true, 1, $false, "11+78=", 2) | Out-Null
$d.Content.Find.Execute("41+43=", $true, $false, $false, $false, $false, $true, 1, $false, "98-15=", 2) | Out-Null
$d.Content.Find.Execute("26+21=", $true, $false, $false, $false, $false, $true, 1, $false, "93-12=", 2) | Out-Null
$d.Content.Find.Execute("8+18=", $true, $false, $false, $false, $false, $true, 1, $false, "93-70=", 2) | Out-Null
$d.Content.Find.Execute("12-0=", $true, $false, $false, $false, $false, $true, 1, $false, "27+10=", 2) | Out-Null
$d.Content.Find.Execute("39+20=", $true, $false, $false, $false, $false, $true, 1, $false, "24+26=", 2) | Out-Null
$d.Content.Find.Execute("64-62=", $true, $false, $false, $false, $false, $true, 1, $false, "73-46=", 2) | Out-Null
$d.Content.Find.Execute("74-43=", $true, $false, $false, $false, $false, $true, 1, $false, "33+44=", 2) | Out-Null
$d.Content.Find.Execute("12+57=", $true, $false, $false, $false, $false, $true, 1, $false, "43-26=", 2) | Out-Null
$d.Content.Find.Execute("76-53=", $true, $false, $false, $false, $false, $true, 1, $false, "59+9=", 2) | Out-Null
$d.Content.Find.Execute("27+38=", $true, $false, $false, $false, $false, $true, 1, $false, "13+51=", 2) | Out-Null
$d.Content.Find.Execute("58-17=", $true, $false, $false, $false, $false, $true, 1, $false, "63+30=", 2) | Out-Null
$d.Content.Find.Execute("61-35=", $true, $false, $false, $false, $false, $true, 1, $false, "68-46=", 2) | Out-Null
$d.Content.Find.Execute("0-0=", $true, $false, $false, $false, $false, $true, 1, $false, "57+7=", 2) | Out-Null
$d.Content.Find.Execute("98-26=", $true, $false, $false, $false, $false, $true, 1, $false, "59-7=", 2) | Out-Null
$d.Content.Find.Execute("72-54=", $true, $false, $false, $false, $false, $true, 1, $false, "22-10=", 2) | Out-Null
$d.Content.Find.Execute("11+0=", $true, $false, $false, $false, $false, $true, 1, $false, "70+11=", 2) | Out-Null
$d.Content.Find.Execute("1+90=", $true, $false, $false, $false, $false, $true, 1, $false, "85-46=", 2) | Out-Null
$d.Content.Find.Execute("43-12=", $true, $false, $false, $false, $false, $true, 1, $false, "80-12=", 2) | Out-Null
$d.Content.Find.Execute("4+28=", $true, $false, $false, $false, $false, $true, 1, $false, "14-10=", 2) | Out-Null
$d.Content.Find.Execute("70-28=", $true, $false, $false, $false, $false, $true, 1, $false, "90-15=", 2) | Out-Null
$d.Content.Find.Execute("4+94=", $true, $false, $false, $false, $false, $true, 1, $false, "71-30=", 2) | Out-Null
$d.Content.Find.Execute("89-9=", $true, $false, $false, $false, $false, $true, 1, $false, "63-19=", 2) | Out-Null
$d.Content.Find.Execute("23+59=", $true, $false, $false, $false, $false, $true, 1, $false, "64+7=", 2) | Out-Null
$d.Content.Find.Execute("84-38=", $true, $false, $false, $false, $false, $true, 1, $false, "8+25=", 2) | Out-Null
$d.Content.Find.Execute("51-5=", $true, $false, $false, $false, $false, $true, 1, $false, "16+14=", 2) | Out-Null
$d.Content.Find.Execute("30+24=", $true, $false, $false, $false, $false, $true, 1, $false, "71+4=", 2) | Out-Null
$d.Content.Find.Execute("47+6=", $true, $false, $false, $false, $false, $true, 1, $false, "61+21=", 2) | Out-Null
$d.Content.Find.Execute("57-33=", $true, $false, $false, $false, $false, $true, 1, $false, "82-15=", 2) | Out-Null
$d.Content.Find.Execute("2+18=", $true, $false, $false, $false, $false, $true, 1, $false, "55+30=", 2) | Out-Null
$d.Content.Find.Execute("29+24=", $true, $false, $false, $false, $false, $true, 1, $false, "60-40=", 2) | Out-Null
$d.Content.Find.Execute("49-29=", $true, $false, $false, $false, $false, $true, 1, $false, "3+70=", 2) | Out-Null
$d.Content.Find.Execute("34+59=", $true, $false, $false, $false, $false, $true, 1, $false, "80-4=", 2) | Out-Null
$d.Content.Find.Execute("16+66=", $true, $false, $false, $false, $false, $true, 1, $false, "96-63=", 2) | Out-Null
$d.Content.Find.Execute("65-54=", $true, $false, $false, $false, $false, $true, 1, $false, "47+44=", 2) | Out-Null
$d.Content.Find.Execute("2+47=", $true, $false, $false, $false, $false, $true, 1, $false, "58-21=", 2) | Out-Null
$d.Content.Find.Execute("90-10=", $true, $false, $false, $false, $false, $true, 1, $false, "0+76=", 2) | Out-Null
$d.Content.Find.Execute("46-28=", $true, $false, $false, $false, $false, $true, 1, $false, "57+35=", 2) | Out-Null
$d.Content.Find.Execute("79-11=", $true, $false, $false, $false, $false, $true, 1, $false, "10+24=", 2) | Out-Null
$d.Content.Find.Execute("84+5=", $true, $false, $false, $false, $false, $true, 1, $false, "9+9=", 2) | Out-Null
$d.Content.Find.Execute("62-44=", $true, $false, $false, $false, $false, $true, 1, $false, "55-30=", 2) | Out-Null
$d.Content.Find.Execute("59-17=", $true, $false, $false, $false, $false, $true, 1, $false, "97-70=", 2) | Out-Null
$d.Content.Find.Execute("51-40=", $true, $false, $false, $false, $false, $true, 1, $false, "65+28=", 2) | Out-Null
$d.Content.Find.Execute("21+69=", $true, $false, $false, $false, $false, $true, 1, $false, "33+53=", 2) | Out-Null
$d.Content.Find.Execute("91-46=", $true, $false, $false, $false, $false, $true, 1, $false, "19+35=", 2) | Out-Null
$d.Content.Find.Execute("94-52=", $true, $false, $false, $false, $false, $true, 1, $false, "45-41=", 2) | Out-Null
$d.Content.Find.Execute("49-9=", $true, $false, $false, $false, $false, $true, 1, $false, "19-10=", 2) | Out-Null
$d.Content.Find.Execute("78-43=", $true, $false, $false, $false, $false, $true, 1, $false, "3+12=", 2) | Out-Null
